$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 already carries the correct border/wrap styling (s=1 for A,B,E; s=2 wrapped for C,D).
# Clone that formatting down onto the two new rows before writing their values.
$ws.Range("A2:E2").Copy()
$ws.Range("A3:E3").PasteSpecial(-4122)
$ws.Range("A4:E4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 2: TC001 becomes the "Validate login functionality" case
$ws.Range("A2").Value = "TC001"
$ws.Range("B2").Value = "Validate login functionality with valid credentials"
$ws.Range("C2").Value = "1. Navigate to the login page`n2. Enter valid username and password`n3.Click the login button`n4.Verify successful login message"
$ws.Range("D2").Value = "User is redirected to the dashboard with a welcome message"
$ws.Range("E2").Value = "As Expected"

# Row 3: the original TC001/CAPTCHA case, renumbered to TC002
$ws.Range("A3").Value = "TC002"
$ws.Range("B3").Value = "Verify that CAPTCHA validation works correctly and prevents automated login attempts"
$ws.Range("C3").Value = "1. Navigate to the login page`n2. Enter valid username and password`n3.Complete the CAPTCHA challenge manually`n4.Click ""Login"" button"
$ws.Range("D3").Value = "User is logged in only if CAPTCHA is solved correctly. If CAPTCHA is incorrect or skipped, login should be blocked with an appropriate error message"
$ws.Range("E3").Value = "As Expected"

# Row 4: new TC003 hardware/LED test case
$ws.Range("A4").Value = "TC003"
$ws.Range("B4").Value = "Verify physical hardware connection and LED status"
$ws.Range("C4").Value = "1. Connect the USB device to the machin`n2. Observe the LED indicator on the device`n3.Confirm the LED blinks three times`n4.Disconnect the device and check for safe removal notification"
$ws.Range("D4").Value = "LED blinks three times and system shows safe removal"
$ws.Range("E4").Value = "As Expected"

$ws.Rows.Item(2).RowHeight = 87
$ws.Rows.Item(3).RowHeight = 87
$ws.Rows.Item(4).RowHeight = 130.5

# Column widths to fit the new, longer content (closest values the width model supports)
$ws.Columns.Item(2).ColumnWidth = 70.5
$ws.Columns.Item(3).ColumnWidth = 22.666666666666664
$ws.Columns.Item(4).ColumnWidth = 39.33333333333333

[void]$ws.Range("D3").Select()
